# Updates cryptos price/volume figures (and swaps the Fetch.AI / Bittensor rows)
# per the Tue Nov  5 18:34:33 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe forces Excel to store a digit-only string (e.g. "569.98")
# as literal text instead of silently coercing it to a Double.

$ws.Range('D2').Value = "69.969.04"
$ws.Range('E2').Value = "  +3.20%  "
$ws.Range('D3').Value = "2.452.63"
$ws.Range('E3').Value = "  +0.92%  "
$ws.Range('E4').Value = "  -0.12%  "
$ws.Range('D5').Value = "'569.98"
$ws.Range('E5').Value = "  +3.14%  "
$ws.Range('D6').Value = "'166.65"
$ws.Range('E6').Value = "  +3.70%  "
$ws.Range('E7').Value = "  -0.16%  "
$ws.Range('D8').Value = "'0.514"
$ws.Range('E8').Value = "  +0.76%  "
$ws.Range('D9').Value = "'0.175"
$ws.Range('E9').Value = "  +11.65%  "
$ws.Range('D10').Value = "2.450.75"
$ws.Range('E10').Value = "  +0.78%  "
$ws.Range('E11').Value = "  -1.65%  "
$ws.Range('E12').Value = "  +2.81%  "
$ws.Range('E13').Value = "  -1.84%  "
$ws.Range('D14').Value = "'0.0000182"
$ws.Range('E14').Value = "  +7.89%  "
$ws.Range('D15').Value = "69.835.54"
$ws.Range('E15').Value = "  +3.12%  "
$ws.Range('D16').Value = "2.902.23"
$ws.Range('E16').Value = "  -0.12%  "
$ws.Range('E17').Value = "  +5.06%  "
$ws.Range('D18').Value = "2.450.71"
$ws.Range('E18').Value = "  +1.02%  "
$ws.Range('D19').Value = "'10.90"
$ws.Range('E19').Value = "  +5.82%  "
$ws.Range('D20').Value = "'7.17"
$ws.Range('E20').Value = "  +5.05%  "
$ws.Range('D21').Value = "'341.63"
$ws.Range('E21').Value = "  +2.12%  "
$ws.Range('E22').Value = "  +3.35%  "
$ws.Range('E23').Value = "  +8.44%  "
$ws.Range('E24').Value = "  -0.12%  "
$ws.Range('D25').Value = "'66.42"
$ws.Range('E25').Value = "  -0.01%  "
$ws.Range('D26').Value = "'3.84"
$ws.Range('E26').Value = "  +5.88%  "
$ws.Range('D27').Value = "2.578.32"
$ws.Range('E27').Value = "  +0.87%  "
$ws.Range('D28').Value = "'8.52"
$ws.Range('E28').Value = "  +5.28%  "
$ws.Range('D29').Value = "'0.984"
$ws.Range('E29').Value = "  -1.54%  "
$ws.Range('D30').Value = "0.0₃0858"
$ws.Range('E30').Value = "  +6.07%  "
$ws.Range('E31').Value = "  +4.06%  "
$ws.Range('B32').Value = "Bittensor"
$ws.Range('C32').Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range('D32').Value = "'461.63"
$ws.Range('E32').Value = "  +9.91%  "
$ws.Range('B33').Value = "Fetch.AI"
$ws.Range('C33').Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range('D33').Value = "'1.25"
$ws.Range('E33').Value = "  +10.26%  "
$ws.Range('E34').Value = "  -0.02%  "
$ws.Range('E35').Value = "  +2.35%  "
$ws.Range('D36').Value = "'160.62"
$ws.Range('E36').Value = "  -0.11%  "
$ws.Range('E37').Value = "  +9.03%  "
$ws.Range('D38').Value = "'19.11"
$ws.Range('E38').Value = "  +0.94%  "
$ws.Range('E39').Value = "  +0.05%  "
$ws.Range('E40').Value = "  +2.31%  "
$ws.Range('E41').Value = "  +3.68%  "
$ws.Range('E42').Value = "  +4.99%  "
$ws.Range('E43').Value = "  +4.13%  "
$ws.Range('D44').Value = "'38.07"
$ws.Range('E44').Value = "  +1.76%  "
$ws.Range('E45').Value = "  +2.53%  "
$ws.Range('E46').Value = "  +6.18%  "
$ws.Range('D47').Value = "'134.30"
$ws.Range('E47').Value = "  +4.12%  "
$ws.Range('E48').Value = "  +1.92%  "
$ws.Range('D49').Value = "'0.0726"
$ws.Range('E49').Value = "  +2.32%  "
$ws.Range('E50').Value = "  +2.75%  "
$ws.Range('D51').Value = "'0.565"
$ws.Range('E51').Value = "  +1.81%  "
